$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Row 2: "goto" destination changes from the old SIMPLIHOME product URL to
#    the Amazon homepage, and becomes a real hyperlink.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "https://www.amazon.com"

# ---------------------------------------------------------------------------
# 2. Remove the bold styling that used to mark the "locatorType" column
#    (C3:C14) -- the bold font is dropped from the workbook entirely.
# ---------------------------------------------------------------------------
$ws.Range("C3:C14").Font.Bold = $false

# ---------------------------------------------------------------------------
# 3. Row 12 used to be a plain click on the "Go" button; it is replaced with
#    a keypress of Enter on the search box instead.
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = "keypress"
$ws.Range("D12").Value = "Search Amazon"
$ws.Range("E12").Value = "input"
$ws.Range("F12").Value = "Enter"

# ---------------------------------------------------------------------------
# 4. Drop the old "#parsersMB2,shipping_amount" assertion row (old row 22) --
#    the firstName/lastName assertions shift up to take rows 22 and 23.
# ---------------------------------------------------------------------------
$ws.Rows("22:22").Delete($xlUp)

# ---------------------------------------------------------------------------
# 5. Append the new "cartassert" rows (24-28) that exercise the cart parser
#    assertions (discounts, sales tax, shipping, cart total, items).
# ---------------------------------------------------------------------------
$newRows = @(
    @{ D = "#parsers,discounts" },
    @{ D = "#parsers,sales_tax" },
    @{ D = "#parsers,shipping_amount" },
    @{ D = "#parsers,cart_total" },
    @{ D = "#parsers,items" }
)

$r = 24
foreach ($row in $newRows) {
    $ws.Range("A$r").Value = "TC001"
    $ws.Range("B$r").Value = "Yes"
    $ws.Range("C$r").Value = "cartassert"
    $ws.Range("D$r").Value = $row.D
    $ws.Range("G$r").Value = 1000
    $ws.Range("H$r").Value = 2000
    $r++
}

# ---------------------------------------------------------------------------
# 6. Rebuild the hyperlinks. The engine does not track hyperlink ranges
#    through row shifts, so clear everything and re-add each link at its
#    final location (F6, F8, D15 stay put; D2 is new; the Weaver mailto link
#    moves from the old F24 up to F23 after the row delete above).
# ---------------------------------------------------------------------------
$ws.Range("A1:H30").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D2"), "https://www.amazon.com")
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:weavernormar@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F8"), "mailto:Welcome@123456")
$ws.Hyperlinks.Add($ws.Range("D15"), "https://www.amazon.com/gp/cart/view.html?ref_=nav_cart")
$ws.Hyperlinks.Add($ws.Range("F23"), "mailto:weavernormar@gmail.com", "", "weavernormar@gmail.com", "Weaver")

# D2 / D15 keep the "Hyperlink" look, matching the rest of the goto links.
$ws.Range("D2").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 7. Selection / scroll position noted in the diff.
# ---------------------------------------------------------------------------
$ws.Range("D16").Select()
$ws.Application.ActiveWindow.ScrollRow = 9

Write-Host "edit complete"
